# Update TPM-derived values in the LR-pairs sheet (Plg-F3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("G2").Value = 0.0345785
$ws.Range("H2").Value = 0.069157
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1564773333333333
$ws.Range("N2").Value = 0.469432
$ws.Range("O2").Value = 0.002686289672660777
$ws.Range("P2").Value = 0.002698380613520143
$ws.Range("Q2").Value = 0.005410751470666666
$ws.Range("R2").Value = 0.032464508824
$ws.Range("S2").Value = 0.002686289672660777
$ws.Range("T2").Value = 0.002698380613520143

# --- Row 3 ---
$ws.Range("G3").Value = 0.0345785
$ws.Range("H3").Value = 0.069157
$ws.Range("O3").Value = 0.9776697736080268
$ws.Range("P3").Value = 0.9820702474411319
$ws.Range("Q3").Value = 1.969232216172833
$ws.Range("R3").Value = 11.815393297037
$ws.Range("S3").Value = 0.9776697736080268
$ws.Range("T3").Value = 0.9820702474411319

# --- Row 4 ---
$ws.Range("G4").Value = 0.0345785
$ws.Range("H4").Value = 0.069157
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.111859
$ws.Range("N4").Value = 0.335577
$ws.Range("O4").Value = 0.001920314400131405
$ws.Range("P4").Value = 0.001928957700248916
$ws.Range("Q4").Value = 0.0038679164315
$ws.Range("R4").Value = 0.023207498589
$ws.Range("S4").Value = 0.001920314400131405
$ws.Range("T4").Value = 0.001928957700248916

# --- Row 5 ---
$ws.Range("G5").Value = 0.0345785
$ws.Range("H5").Value = 0.069157
$ws.Range("M5").Value = 0.7830269999999999
$ws.Range("N5").Value = 1.566054
$ws.Range("O5").Value = 0.0134424411427931
$ws.Range("P5").Value = 0.009001963550260048
$ws.Range("Q5").Value = 0.0270758991195
$ws.Range("R5").Value = 0.108303596478
$ws.Range("S5").Value = 0.0134424411427931
$ws.Range("T5").Value = 0.009001963550260048

# --- Row 6 ---
$ws.Range("G6").Value = 0.0345785
$ws.Range("H6").Value = 0.069157
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.2493803333333333
$ws.Range("N6").Value = 0.7481409999999999
$ws.Range("O6").Value = 0.004281181176387862
$ws.Range("P6").Value = 0.004300450694838813
$ws.Range("Q6").Value = 0.008623197856166666
$ws.Range("R6").Value = 0.05173918713699999
$ws.Range("S6").Value = 0.004281181176387862
$ws.Range("T6").Value = 0.004300450694838813
